# Applies the "Update docx golden tests for style changes" edit:
#   1. Adds a new paragraph style "AbstractTitle" ("Abstract Title") that
#      sits between the existing "Date" and "Abstract" styles.
#   2. Tightens the space-before on the "Abstract" style from 15pt (300)
#      to 5pt (100), leaving space-after untouched.
#   3. Gives the "ImportTok" character style a green, bold run format and
#      the "BuiltInTok" character style a green run format.

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" paragraph style -----------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. "Abstract" style: space-before 300 -> 100 ------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. Pandoc syntax-highlighting token styles --------------------------
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768
